$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the third duplicate "person_id_nbr" block (rows 8-10), which shifts all
# subsequent "year" rows up by 3 and shrinks the used range from A1:C636 to A1:C633.
$ws.Rows("8:10").Delete()
